$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rewrite rows 2-19: 3 brand-new scraped listings prepended, the former
# rows 2-16 shifted down by 3 and re-stamped with the new scrape time.
$ws.Cells.Item(2, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(2, 2).Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Cells.Item(2, 3).Value = 'システム開発'
$ws.Cells.Item(2, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = '期限情報なし'
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Cells.Item(2, 7).Value = 243
$ws.Cells.Item(2, 8).Value = '🔥API ◆ツール'

$ws.Cells.Item(3, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(3, 2).Value = '【高単価業務自動化】行政書士向けシステム開発依頼'
$ws.Cells.Item(3, 3).Value = 'システム開発'
$ws.Cells.Item(3, 4).Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = '期限情報なし'
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5441252'
$ws.Cells.Item(3, 7).Value = 210
$ws.Cells.Item(3, 8).Value = '◆開発,システム開発'

$ws.Cells.Item(4, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(4, 2).Value = '初回 FastAPIバックエンドの軽微な修正・調整対応エンジニア募集'
$ws.Cells.Item(4, 3).Value = 'システム開発'
$ws.Cells.Item(4, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = '期限情報なし'
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5441207'
$ws.Cells.Item(4, 7).Value = 183
$ws.Cells.Item(4, 8).Value = '🔥API'

$ws.Cells.Item(5, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(5, 2).Value = '【募集】RPAツール「RoboTANGO」設定代行の専門家を探しています'
$ws.Cells.Item(5, 3).Value = 'システム開発'
$ws.Cells.Item(5, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = '期限情報なし'
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5405023'
$ws.Cells.Item(5, 7).Value = 178
$ws.Cells.Item(5, 8).Value = '★bot ◆ツール'

$ws.Cells.Item(6, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(6, 2).Value = '【急募】縫製工場向けPL・CF可視化アプリのMVP開発'
$ws.Cells.Item(6, 3).Value = 'システム開発'
$ws.Cells.Item(6, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = '期限情報なし'
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5440957'
$ws.Cells.Item(6, 7).Value = 93
$ws.Cells.Item(6, 8).Value = '◆開発 ◇アプリ'

$ws.Cells.Item(7, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(7, 2).Value = '【長期募集】クラウドサービス開発・保守エンジニアを求む!'
$ws.Cells.Item(7, 3).Value = 'システム開発'
$ws.Cells.Item(7, 4).Value = '1,000 ~ 5,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = '期限情報なし'
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5440461'
$ws.Cells.Item(7, 7).Value = 60
$ws.Cells.Item(7, 8).Value = '◆開発'

$ws.Cells.Item(8, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(8, 2).Value = '【急募】古いPHPとPerlプログラムのアップデート依頼'
$ws.Cells.Item(8, 3).Value = 'システム開発'
$ws.Cells.Item(8, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = '期限情報なし'
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5440861'
$ws.Cells.Item(8, 7).Value = 33
$ws.Cells.Item(8, 8).Value = '○PHP'

$ws.Cells.Item(9, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(9, 2).Value = '【急募】シティヘブンの出勤情報を自動取得・管理したい!'
$ws.Cells.Item(9, 3).Value = 'システム開発'
$ws.Cells.Item(9, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = '期限情報なし'
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5440436'
$ws.Cells.Item(9, 7).Value = 33
$ws.Cells.Item(9, 8).Value = '◇管理'

$ws.Cells.Item(10, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(10, 2).Value = '簡易サイト修正・その他小規模タスク依頼'
$ws.Cells.Item(10, 3).Value = 'システム開発'
$ws.Cells.Item(10, 4).Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = '期限情報なし'
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5441146'
$ws.Cells.Item(10, 7).Value = 30
$ws.Cells.Item(10, 8).Value = '◇サイト'

$ws.Cells.Item(11, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(11, 2).Value = '進行管理およびチームディレクションを担当'
$ws.Cells.Item(11, 3).Value = 'システム開発'
$ws.Cells.Item(11, 4).Value = '~ 5,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = '期限情報なし'
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Cells.Item(11, 7).Value = 30
$ws.Cells.Item(11, 8).Value = '◇管理'

$ws.Cells.Item(12, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(12, 2).Value = '初回 n8n+Gemini+Typefully+GoogleスプレッドのX/Threads自動投稿システム'
$ws.Cells.Item(12, 3).Value = 'システム開発'
$ws.Cells.Item(12, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = '期限情報なし'
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5440440'
$ws.Cells.Item(12, 7).Value = 33

$ws.Cells.Item(13, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(13, 2).Value = '急募 限定公開 PR 限定公開の仕事'
$ws.Cells.Item(13, 3).Value = 'システム開発'
$ws.Cells.Item(13, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = '期限情報なし'
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5440230'
$ws.Cells.Item(13, 7).Value = 25

$ws.Cells.Item(14, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(14, 2).Value = '〖リモート可〗Delphiエンジニア募集'
$ws.Cells.Item(14, 3).Value = 'システム開発'
$ws.Cells.Item(14, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = '期限情報なし'
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5341051'
$ws.Cells.Item(14, 7).Value = 25

$ws.Cells.Item(15, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(15, 2).Value = '【リーダー募集×リモートOK】SRE/インフラエンジニア(Google Cloud/長期金融系案件)'
$ws.Cells.Item(15, 3).Value = 'システム開発'
$ws.Cells.Item(15, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = '期限情報なし'
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5435080'
$ws.Cells.Item(15, 7).Value = 25

$ws.Cells.Item(16, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(16, 2).Value = '【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)'
$ws.Cells.Item(16, 3).Value = 'システム開発'
$ws.Cells.Item(16, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(16, 5).Value = '期限情報なし'
$ws.Cells.Item(16, 6).Value = 'https://www.lancers.jp/work/detail/5435079'
$ws.Cells.Item(16, 7).Value = 25

$ws.Cells.Item(17, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(17, 2).Value = '当組織のエンジニア追加募集。'
$ws.Cells.Item(17, 3).Value = 'システム開発'
$ws.Cells.Item(17, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(17, 5).Value = '期限情報なし'
$ws.Cells.Item(17, 6).Value = 'https://www.lancers.jp/work/detail/5441084'
$ws.Cells.Item(17, 7).Value = 18

$ws.Cells.Item(18, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(18, 2).Value = '【急募】弊社Websiteの保守運用をお任せできる方を探しています!'
$ws.Cells.Item(18, 3).Value = 'システム開発'
$ws.Cells.Item(18, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(18, 5).Value = '期限情報なし'
$ws.Cells.Item(18, 6).Value = 'https://www.lancers.jp/work/detail/5440806'
$ws.Cells.Item(18, 7).Value = 18

$ws.Cells.Item(19, 1).Value = '2025-11-26 01:20:19'
$ws.Cells.Item(19, 2).Value = '【急募】n8n ワークフロー実装とGoogle Sheets作成(70万円 ~ )'
$ws.Cells.Item(19, 3).Value = 'システム開発'
$ws.Cells.Item(19, 4).Value = '1,000 ~ 5,000 円 / 固定'
$ws.Cells.Item(19, 5).Value = '期限情報なし'
$ws.Cells.Item(19, 6).Value = 'https://www.lancers.jp/work/detail/5441082'
$ws.Cells.Item(19, 7).Value = 10

# Rows 17-19 are brand new beyond the former A1:H16 extent, so the F column
# there needs the Hyperlink cell style explicitly (F2:F16 already carry it).
$ws.Range("F17").Style = "Hyperlink"
$ws.Range("F18").Style = "Hyperlink"
$ws.Range("F19").Style = "Hyperlink"

# Hyperlinks cannot be edited/removed individually in this engine, so rebuild
# the whole collection from scratch in the final top-to-bottom order.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5441252')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5441207')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5405023')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5440957')
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5440461')
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5440861')
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5440436')
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5441146')
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5418064')
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5440440')
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5440230')
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5341051')
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5435080')
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5435079')
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5441084')
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5440806')
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5441082')

# Column D widened from 30 to 32 characters. ColumnWidth pads by ~5/6 of a
# character over the stored <col width>, so compensate to land exactly on 32.
$ws.Columns.Item(4).ColumnWidth = (32 - 0.8333333333333334)